$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.907.95"
$ws.Range("E2").Value = "  +2.50%  "

$ws.Range("D3").Value = "1.875.03"
$ws.Range("E3").Value = "  +0.82%  "

$ws.Range("E4").Value = "  -0.81%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4849"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3808"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07372"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.67%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9401"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.96"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.11%  "

$ws.Range("E12").Value = "  -1.01%  "

$ws.Range("D13").Value = "1.913.23"
$ws.Range("E13").Value = "  +2.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.545"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.595"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.91%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.014"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008873"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.63%  "

$ws.Range("E19").Value = "  -0.77%  "

$ws.Range("D20").Value = "27.927.35"
$ws.Range("E20").Value = "  +2.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.73%  "

$ws.Range("E22").Value = "  +0.29%  "

$ws.Range("D23").Value = "2.122.23"
$ws.Range("E23").Value = "  +1.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.41%  "

$ws.Range("E26").Value = "  -0.98%  "

$ws.Range("E27").Value = "  +0.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.971"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08899"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.337"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.225"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7719"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.643"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.30%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.723"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.37%  "

$ws.Range("E37").Value = "  +0.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02046"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5594"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05371"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.004"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.051"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.529"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1519"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.62%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4886"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.84%  "

$ws.Range("E48").Value = "  -0.87%  "

$ws.Range("E49").Value = "  +2.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06114"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.61%  "
